# Weekly driver report update for 2025-04-21
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: MediaTek ... 3.3.0.800 -- Good Roaming % updated
$ws.Range("D3").Value = 89.40000000000001

# Row 4: MediaTek ... 3.3.0.897 -- Good Roaming % updated
$ws.Range("D4").Value = 93.7

# Row 5: now MediaTek ... 3.3.0.824 (was Intel AX211 23.60.1.2)
$ws.Range("A5").Value = "MediaTek Wi-Fi 6E MT7922 (RZ616) 160MHz PCIe Adapter - 3.3.0.824"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 29
$ws.Range("D5").Value = 97.59999999999999

# Row 6: now Intel AX211 23.60.1.2 (was MediaTek ... 3.3.0.824)
$ws.Range("A6").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.60.1.2"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 18
$ws.Range("D6").Value = 97.7

# Row 7: now Intel AX201 23.40.0.4 (was MediaTek ... 3.3.0.908)
$ws.Range("A7").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.40.0.4"
$ws.Range("B7").Value = 35
$ws.Range("C7").Value = 511
$ws.Range("D7").Value = 98.40000000000001

# Row 8: now MediaTek ... 3.3.0.908 (was Intel AX201 23.40.0.4)
$ws.Range("A8").Value = "MediaTek Wi-Fi 6E MT7922 (RZ616) 160MHz PCIe Adapter - 3.3.0.908"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 35
$ws.Range("D8").Value = 98.5

# Row 9: Totals updated
$ws.Range("B9").Value = 50
$ws.Range("C9").Value = 979
